# Enter the 2025-05-13 batch of oyster size measurements under the
# existing 2025-05-06 batch (rows 27-51), reusing the same five-family
# sample-naming scheme (A1-A5, B1-B5, C1-C5, D1-D5, E1-E5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$samples = @("A1","A2","A3","A4","A5","B1","B2","B3","B4","B5","C1","C2","C3","C4","C5","D1","D2","D3","D4","D5","E1","E2","E3","E4","E5")

$lengths = @(
    17.48,
    14.289,
    15.5,
    14.762,
    16.318999999999999,
    18.527000000000001,
    15.419,
    14.48,
    16.559000000000001,
    15.683999999999999,
    20.678000000000001,
    14.311,
    14.215,
    17.335000000000001,
    21.035,
    19.798999999999999,
    17.22,
    16.512,
    16.361000000000001,
    21.093,
    19.13,
    15.736000000000001,
    16.893000000000001,
    15.061,
    16.948
)

$areas = @(
    177.47900000000001,
    113.44199999999999,
    123.788,
    107.637,
    184.042,
    169.06399999999999,
    140.43,
    143.61099999999999,
    154.05500000000001,
    174.322,
    237.83600000000001,
    127.06,
    140.88499999999999,
    194.83099999999999,
    258.553,
    219.57900000000001,
    160.6,
    175.94300000000001,
    147.91,
    234.952,
    166.50700000000001,
    158.62,
    131.79,
    88.245000000000005,
    137.267
)

# Rows 27-51 already carried an empty, styled "B" cell left over from the
# sheet's original blank template; clear that leftover style so the newly
# entered sample-name cells land unstyled.
$ws.Range("B27:B51").ClearFormats()

$startRow = 27
for ($i = 0; $i -lt $samples.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = 20250513
    $ws.Cells.Item($r, 2).Value = $samples[$i]
    $ws.Cells.Item($r, 3).Value = $lengths[$i]
    $ws.Cells.Item($r, 4).Value = $areas[$i]
}

# Leave the selection where data entry ended.
$ws.Range("B30").Select() | Out-Null
